$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as TEXT (matches original inline-string cells),
# even when the string looks like a number (e.g. "1.002"), then restore default styling.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Rows 2-31: price (D) and volume (E) updates only
Set-TextValue $ws.Range("D2") '24.724.45'
Set-TextValue $ws.Range("E2") '  +0.87%  '
Set-TextValue $ws.Range("D3") '1.694.47'
Set-TextValue $ws.Range("E3") '  -0.01%  '
Set-TextValue $ws.Range("D4") '1.002'
Set-TextValue $ws.Range("E4") '  -0.03%  '
Set-TextValue $ws.Range("D5") '317.19'
Set-TextValue $ws.Range("E5") '  +1.19%  '
Set-TextValue $ws.Range("D6") '1.003'
Set-TextValue $ws.Range("E6") '  +0.05%  '
Set-TextValue $ws.Range("D7") '0.3951'
Set-TextValue $ws.Range("E7") '  +0.01%  '
Set-TextValue $ws.Range("D8") '0.4069'
Set-TextValue $ws.Range("E8") '  +0.74%  '
Set-TextValue $ws.Range("D9") '1.492'
Set-TextValue $ws.Range("E9") '  -2.00%  '
Set-TextValue $ws.Range("D10") '1.001'
Set-TextValue $ws.Range("E10") '  -0.11%  '
Set-TextValue $ws.Range("D11") '52.17'
Set-TextValue $ws.Range("E11") '  -2.62%  '
Set-TextValue $ws.Range("D12") '0.08903'
Set-TextValue $ws.Range("E12") '  +1.40%  '
Set-TextValue $ws.Range("D13") '7.270'
Set-TextValue $ws.Range("E13") '  -0.46%  '
Set-TextValue $ws.Range("D14") '23.70'
Set-TextValue $ws.Range("E14") '  +2.32%  '
Set-TextValue $ws.Range("D15") '8.075'
Set-TextValue $ws.Range("E15") '  +7.25%  '
Set-TextValue $ws.Range("D16") '0.00001322'
Set-TextValue $ws.Range("E16") '  +0.14%  '
Set-TextValue $ws.Range("D17") '1.697.13'
Set-TextValue $ws.Range("E17") '  +0.27%  '
Set-TextValue $ws.Range("D18") '99.88'
Set-TextValue $ws.Range("E18") '  -0.42%  '
Set-TextValue $ws.Range("D19") '0.07034'
Set-TextValue $ws.Range("E19") '  -0.89%  '
Set-TextValue $ws.Range("D20") '19.64'
Set-TextValue $ws.Range("E20") '  +1.08%  '
Set-TextValue $ws.Range("D21") '7.015'
Set-TextValue $ws.Range("E21") '  +4.78%  '
Set-TextValue $ws.Range("D22") '1.007'
Set-TextValue $ws.Range("E22") '  +0.39%  '
Set-TextValue $ws.Range("D23") '14.40'
Set-TextValue $ws.Range("E23") '  +1.78%  '
Set-TextValue $ws.Range("D24") '24.697.43'
Set-TextValue $ws.Range("E24") '  +0.84%  '
Set-TextValue $ws.Range("D25") '3.215'
Set-TextValue $ws.Range("E25") '  +6.93%  '
Set-TextValue $ws.Range("D26") '2.361'
Set-TextValue $ws.Range("E26") '  +1.80%  '
Set-TextValue $ws.Range("D27") '22.81'
Set-TextValue $ws.Range("E27") '  +1.93%  '
Set-TextValue $ws.Range("D28") '162.39'
Set-TextValue $ws.Range("E28") '  +1.96%  '
Set-TextValue $ws.Range("D29") '136.13'
Set-TextValue $ws.Range("E29") '  +1.76%  '
Set-TextValue $ws.Range("D30") '5.186'
Set-TextValue $ws.Range("E30") '  +0.25%  '
Set-TextValue $ws.Range("D31") '7.588'
Set-TextValue $ws.Range("E31") '  +1.09%  '

# Rows 32-51: coin list shifted up by one (WrappedliquidstakedEther2.0 dropped out of
# top 50, Cronos entered at the bottom) -- update B (name), C (link), D (price), E (volume)
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D32") '0.08631'
Set-TextValue $ws.Range("E32") '  -0.22%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D33") '1.059'
Set-TextValue $ws.Range("E33") '  -2.57%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D34") '7.106'
Set-TextValue $ws.Range("E34") '  -3.19%  '
$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D35") '11.36'
Set-TextValue $ws.Range("E35") '  +2.95%  '
$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D36") '0.2738'
Set-TextValue $ws.Range("E36") '  +0.88%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D37") '14.53'
Set-TextValue $ws.Range("E37") '  -1.45%  '
$ws.Range("B38").Value = 'WEMIXTOKEN'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D38") '1.885'
Set-TextValue $ws.Range("E38") '  -4.40%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D39") '0.09219'
Set-TextValue $ws.Range("E39") '  +2.54%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D40") '0.02729'
Set-TextValue $ws.Range("E40") '  -1.14%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D41") '1.474'
Set-TextValue $ws.Range("E41") '  -0.14%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D42") '0.7690'
Set-TextValue $ws.Range("E42") '  +0.60%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D43") '16.07'
Set-TextValue $ws.Range("E43") '  +3.53%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D44") '2.610'
Set-TextValue $ws.Range("E44") '  +6.66%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range("D45") '0.7170'
Set-TextValue $ws.Range("E45") '  +0.24%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D46") '4.229'
Set-TextValue $ws.Range("E46") '  +1.51%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D47") '1.002'
Set-TextValue $ws.Range("E47") '  +0.01%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D48") '140.24'
Set-TextValue $ws.Range("E48") '  +0.10%  '
$ws.Range("B49").Value = 'Flow'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
Set-TextValue $ws.Range("D49") '1.323'
Set-TextValue $ws.Range("E49") '  +1.92%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D50") '90.81'
Set-TextValue $ws.Range("E50") '  +5.17%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D51") '0.07979'
Set-TextValue $ws.Range("E51") '  -0.12%  '
